$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.178.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.49%  '
$ws.Range("D3").Value = '''1.923.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.02%  '
$ws.Range("D4").Value = '''0.9986'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''245.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.87%  '
$ws.Range("D6").Value = '''0.7206'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.17%  '
$ws.Range("D7").Value = '''0.9982'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '''0.3247'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.79%  '
$ws.Range("D9").Value = '''26.44'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.84%  '
$ws.Range("D10").Value = '''0.06841'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").Value = '''0.7952'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.52%  '
$ws.Range("D12").Value = '''0.07922'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '''1.920.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.17%  '
$ws.Range("D14").Value = '''5.394'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").Value = '''94.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.63%  '
$ws.Range("D16").Value = '''14.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("D17").Value = '''259.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.16%  '
$ws.Range("D18").Value = '''30.180.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.43%  '
$ws.Range("D19").Value = '''5.834'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").Value = '''0.000007950'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = '''2.170.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.20%  '
$ws.Range("D22").Value = '''0.9986'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = '''0.9980'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = '''6.860'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("D25").Value = '''9.676'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").Value = '''160.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.19%  '
$ws.Range("D27").Value = '''0.1338'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.77%  '
$ws.Range("D28").Value = '''18.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.53%  '
$ws.Range("D29").Value = '''2.247'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.60%  '
$ws.Range("D30").Value = '''1.353'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = '''1.545'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.27%  '
$ws.Range("D32").Value = '''4.415'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.05%  '
$ws.Range("D33").Value = '''4.181'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.40%  '
$ws.Range("D34").Value = '''0.05040'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("E35").Value = '  -1.76%  '
$ws.Range("D36").Value = '''0.7361'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("D37").Value = '''2.731'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("D38").Value = '''0.01938'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("D39").Value = '''2.805'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").Value = '''80.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.50%  '
$ws.Range("D41").Value = '''6.500'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("D42").Value = '''0.4429'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.63%  '
$ws.Range("D43").Value = '''2.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("D44").Value = '''0.9986'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '''0.8297'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("D46").Value = '''102.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("D47").Value = '''9.707'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.02%  '
$ws.Range("D48").Value = '''7.260'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.82%  '
$ws.Range("D49").Value = '''36.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").Value = '''0.4097'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.36%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.473'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.59%  '
